$wb = $excel.ActiveWorkbook

# Helper: force a range to hold a *text* value (even if it looks numeric,
# e.g. "001917" or "4.16") without leaving a lasting custom number-format
# style behind - we briefly flip on a text format, assign, then put the
# cell style back to the workbook default ("Normal").
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating the "2022-Q2" sheet
#    (same layout/header/styles) and placing it right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# Overwrite the data rows of the new sheet with the 2022-Q3 fund data.
Set-TextValue $q3.Range("B2") "001917"
$q3.Range("C2").Value = "招商量化精选股票A"
Set-TextValue $q3.Range("D2") "4.16"
Set-TextValue $q3.Range("E2") "92.70"
Set-TextValue $q3.Range("F2") "1.37"
Set-TextValue $q3.Range("G2") "0.0570"
$q3.Range("H2").Value = 10

Set-TextValue $q3.Range("B3") "007950"
$q3.Range("C3").Value = "招商量化精选股票C"
Set-TextValue $q3.Range("D3") "2.39"
Set-TextValue $q3.Range("E3") "92.70"
Set-TextValue $q3.Range("F3") "1.37"
Set-TextValue $q3.Range("G3") "0.0327"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Insert a new row into the "总计" (Total) summary sheet for 2022-Q3,
#    right after the header row, pushing all the other rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Restore the formatting on the new row 2 (copy from row 3, which still
# carries the original - pre-insert - styling for each column).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.09

# Renumber the index column (A) for the rows that shifted down.
for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally-active tab (copying the sheet made the new copy
# active; put the selection back on "2020-Q4", which was active before).
$wb.Worksheets.Item("2020-Q4").Activate()

